# 989: Adds CMS to test imports and cms table to extract process
#
# Adds a new "CMS" worksheet at the end of the workbook, populates the
# header row with the CMS extract column names, formats that header row
# (Arial 9pt / #333333 font on a white fill, left aligned, taller row),
# and leaves the new sheet selected/active - mirroring a normal
# "add a worksheet for CMS data" edit made interactively in Excel.

$wb = $excel.ActiveWorkbook

# --- add the new "CMS" sheet as the last tab in the workbook ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cms = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$cms.Name = "CMS"

# --- header row values -------------------------------------------------------
$headers = @(
    "Contact_ID`n",
    "Contact_Date`n",
    "Contact_Type_Code",
    "Contact_Type_Desc",
    "Contact_Staff_Name",
    "Contact_Staff_Key",
    "Contact_Staff_Grade",
    "Contact_Team_Key",
    "Contact_Provider_Code",
    "OM_Name`n",
    "OM_Key`n",
    "OM_Grade`n",
    "OM_Team_Key`n",
    "OM_Provider_Code`n"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cms.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- format the header row ---------------------------------------------------
# Build the look (Arial 9, color #333333, solid white fill, left aligned) on
# a scratch cell once, then stamp that exact format onto the header range in
# a single paste so every header cell ends up sharing one cell style.
$helper = $cms.Range("ZZ1")
$helper.Font.Name = "Arial"
$helper.Font.Size = 9
$helper.Font.Color = 3355443
$helper.Interior.Color = 16777215
$helper.Interior.PatternColor = 16777215
$helper.HorizontalAlignment = -4131

$headerRange = $cms.Range("A1:N1")
$helper.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null
$helper.Clear() | Out-Null

$cms.Rows.Item(1).RowHeight = 23.25

# --- leave selection similar to how the sheet was left in the workbook -----
$cms.Range("D32").Select() | Out-Null
